$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores numeric-looking values (cpu/mem utilization, vc count)
# as plain text. Assigning a numeric-looking string straight to .Value
# lets Excel auto-convert it to a real number, so force Text format
# before writing, then restore the cell's original style/format
# afterwards so no visible formatting changes stick around.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# F2: master RE cpu utilization 40 -> 28
Set-TextValue $ws.Range("F2") "28"

# Row 4: F4 (cpu) 31 -> 61, G4 (memory) 38 -> 39
Set-TextValue $ws.Range("F4") "61"
Set-TextValue $ws.Range("G4") "39"

# New row 5 - newly supported ex3400 switch
$ws.Range("A5").Value = "oren-flr1sw-B1"
$ws.Range("B5").Value = "JUNOS OS Kernel 32-bit  [20180119.e26d166_builder_master]"
$ws.Range("C5").Value = "ex3400-48p"
$ws.Range("D5").Value = "10.9.106.11/23"
Set-TextValue $ws.Range("E5") "4"
Set-TextValue $ws.Range("F5") "23"
Set-TextValue $ws.Range("G5") "18"

# Selection now spans the freshly-added data
$ws.Range("A2:G14").Select()

# Reposition the application window (workbookView xWindow/yWindow)
$win = $excel.ActiveWindow
$win.Left = 5010
$win.Top = 2730
